# Generate Report for Handoff
#
# The "cad4860f-87f8-48a8-8617-6517a6f51269.md" entry (row 7 on every
# sheet) moved from "In Translation" to "Ready for handoff" and its
# handoff timestamps were refreshed.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E7").Value = "Ready for handoff"
$wsOverview.Range("F7").Value = "Ready for handoff"
$wsOverview.Range("G7").Value = "2016-08-22 00:55:50"

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C7").Value = "Ready for handoff"
$wsZhCn.Range("H7").Value = "2016-08-22 00:55:45"

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C7").Value = "Ready for handoff"
$wsDeDe.Range("H7").Value = "2016-08-22 00:55:50"
